$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in values for row 9 (Жуков Никита), columns C:G with 5
$ws.Range("C9:G9").Value = 5

# Update active cell/selection to H12
$ws.Range("H12").Select()
